$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Remove the standalone "Meta description" paragraph that currently
#    sits right under the H1 title ("Play Dynamite Frenzy Free & Review
#    | RTP 95.95%").
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Meta description")) {
        $p.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------
# 2. Insert a new bold "title" paragraph right before the final
#    paragraph (the one that used to hold the image-generation
#    "Prompt: ..." text).
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
$last = $d.Paragraphs.Item($count)
$last.Range.InsertParagraphBefore()

$count = $d.Paragraphs.Count
$newp = $d.Paragraphs.Item($count - 1)

$titleXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr>' +
    '<w:t>Play Dynamite Frenzy Free &amp; Review | RTP 95.95%</w:t></w:r></w:p></w:body>' +
    '</w:document></pkg:xmlData></pkg:part></pkg:package>'

$newp.Range.InsertXML($titleXml)

# ---------------------------------------------------------------------
# 3. Swap out the old AI-image "Prompt: ..." copy for the real meta
#    description text, keeping the existing italic run formatting.
# ---------------------------------------------------------------------
$oldPrompt = 'Prompt: Please create a feature image for Dynamite Frenzy that includes a happy Maya warrior with glasses in a cartoon style. Description: The feature image should be vibrant and eye-catching. It should include the Maya warrior prominently in the center, happily holding a stack of golden coins with a big smile on their face. The Maya warrior should be wearing glasses, which should sparkle in the light. There should be colorful gems and gold nuggets scattered around the Maya warrior, giving the impression that they are in a mine. The background should be dark, and the overall image should be in a cartoon style, with bold outlines and bright colors. The Maya warrior should be surrounded by elements from the game, such as dynamite sticks, pickaxes, and lanterns. The title of the game, "Dynamite Frenzy", should be prominently displayed above the image in bold, white letters with an explosion effect.'
$newDescription = 'Read our in-depth review of Dynamite Frenzy, the online slot game featuring a unique mining theme and a Free Spin feature. Play for free and win with an RTP of 95.95%.'

$d.Content.Find.Execute($oldPrompt, $true, $false, $false, $false, $false,
                         $true, 1, $false, $newDescription, 2)
